$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $text) {
    $r = $ws.Range($range)
    $r.NumberFormat = "@"
    $r.Value = $text
    $r.ClearFormats()
}

# Update N values (904 -> 880) for the six existing variable rows
Set-TextValue "B4" "880"
Set-TextValue "B6" "880"
Set-TextValue "B8" "880"
Set-TextValue "B10" "880"
Set-TextValue "B12" "880"
Set-TextValue "B14" "880"

# Add new row 16 data for the "ingreso" variable
$ws.Range("A16").Value = "ingreso"
Set-TextValue "B16" "880"
Set-TextValue "C16" "2.277"
Set-TextValue "D16" "16"
Set-TextValue "E16" "3.102"
Set-TextValue "F16" "896"
$ws.Range("G16").Value = "0.825**"

# Add new row 17 (standard errors for "ingreso"); the other cells in this row
# are explicit (empty) text cells, matching the rest of the table's layout.
$ws.Range("A17").Value = ""
$ws.Range("B17").Value = ""
$ws.Range("C17").Value = "(0.044)"
$ws.Range("D17").Value = ""
$ws.Range("E17").Value = "(0.103)"
$ws.Range("F17").Value = ""
$ws.Range("G17").Value = ""

# Move the footnote text down to row 18, updating it to include "ingreso"
$ws.Range("A18").Value = "If the table includes missing values (.n, .o, .v etc.) see the Missing values section in the help file for the Stata command iebaltab for definitions of these values. Significance: ***=.01, **=.05, *=.1. Full user input as written by user: [iebaltab dummy_jb dummy_d1 dummy_ara cantidad_jb cantidad_d1 cantidad_ara ingreso , groupvar(dummy_oxxo) control(0) savexlsx(difmedias_controles_staggered_variables_2011) replace] "
$ws.Range("B18").Value = ""
$ws.Range("C18").Value = ""
$ws.Range("D18").Value = ""
$ws.Range("E18").Value = ""
$ws.Range("F18").Value = ""
$ws.Range("G18").Value = ""
